# Apply the target edit:
#  - Row 3: clear G3 (was "yes"), rename H3's shared string from 洪國瑋 to 鄭守利
#  - Row 4 (new): A4 = victor.hou@kingza.com.tw (hyperlink, same style as A2/A3),
#                 B4 = "yes", C4 = "yes", G4 = "yes"
#  - Selection moves to G4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename H3's shared string value ("洪國瑋" -> "鄭守利")
$ws.Range("H3").Value = "鄭守利"

# Remove the stray "yes" that was in G3
$ws.Range("G3").ClearContents()

# Build the new row 4
$ws.Range("A4").Value = "victor.hou@kingza.com.tw"
$ws.Range("B4").Value = "yes"
$ws.Range("C4").Value = "yes"
$ws.Range("G4").Value = "yes"

# Hyperlink A4 the same way A2/A3 already are, then restore the
# "superlink" style (Hyperlinks.Add re-applies its own style, so set
# it again afterwards to match A3 exactly).
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:victor.hou@kingza.com.tw")
$ws.Range("A4").Style = $ws.Range("A3").Style

# Update selection to match the saved workbook state
$ws.Range("G4").Select() | Out-Null
